$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.506.53"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -6.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.557.14"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "396.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.546.94"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -9.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.680"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -11.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -18.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000325"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -19.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.93"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -8.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.122.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.18"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -7.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.136"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.89%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.82"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +7.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.536.91"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.69"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -8.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.592.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.02"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -10.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "394.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -12.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.89"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -9.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.92"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.92"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.42"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.95"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -11.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.01"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.73"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -16.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.59"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.112"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.79"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.148"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.62"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.89"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0436"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -11.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0653"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -12.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -12.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.130"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -11.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.08"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +16.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "141.69"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.94"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.96"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.07"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -9.63%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.48"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.58%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.06"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.65"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.276"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -9.66%  "
